$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.478.57'
$ws.Range('E2').Value = '  -2.98%  '
$ws.Range('D3').Value = '1.992.69'
$ws.Range('E3').Value = '  -4.92%  '
$ws.Range('D4').Value = '1.013'
$ws.Range('E4').Value = '  +0.95%  '
$ws.Range('D5').Value = '329.09'
$ws.Range('E5').Value = '  -4.09%  '
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('D7').Value = '0.5012'
$ws.Range('E7').Value = '  -4.52%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4230'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -4.41%  '
$ws.Range('D9').Value = '53.81'
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('D10').Value = '0.08913'
$ws.Range('E10').Value = '  -4.70%  '
$ws.Range('E11').Value = '  -5.04%  '
$ws.Range('D12').Value = '23.15'
$ws.Range('E12').Value = '  -6.43%  '
$ws.Range('D13').Value = '2.007.63'
$ws.Range('E13').Value = '  -6.48%  '
$ws.Range('D14').Value = '7.962'
$ws.Range('E14').Value = '  -7.15%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.450'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -6.84%  '
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = '94.17'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001111'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -4.22%  '
$ws.Range('D19').Value = '0.06758'
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('D20').Value = '19.36'
$ws.Range('E20').Value = '  -8.64%  '
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('D22').Value = '5.927'
$ws.Range('E22').Value = '  -6.44%  '
$ws.Range('D23').Value = '29.505.59'
$ws.Range('E23').Value = '  -3.04%  '
$ws.Range('E24').Value = '  -3.85%  '
$ws.Range('D25').Value = '2.323'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D27').Value = '156.89'
$ws.Range('E27').Value = '  -3.89%  '
$ws.Range('D28').Value = '6.271'
$ws.Range('E28').Value = '  -7.51%  '
$ws.Range('D29').Value = '2.299'
$ws.Range('E29').Value = '  -8.50%  '
$ws.Range('D30').Value = '127.63'
$ws.Range('E30').Value = '  -4.43%  '
$ws.Range('D31').Value = '1.058'
$ws.Range('E31').Value = '  -7.00%  '
$ws.Range('D32').Value = '0.09923'
$ws.Range('E32').Value = '  -5.45%  '
$ws.Range('D33').Value = '1.545'
$ws.Range('E33').Value = '  -6.20%  '
$ws.Range('D34').Value = '5.823'
$ws.Range('E34').Value = '  -7.21%  '
$ws.Range('E35').Value = '  -2.05%  '
$ws.Range('D36').Value = '0.02458'
$ws.Range('E36').Value = '  -6.73%  '
$ws.Range('D37').Value = '9.194'
$ws.Range('E37').Value = '  -9.67%  '
$ws.Range('D38').Value = '0.06375'
$ws.Range('E38').Value = '  -6.40%  '
$ws.Range('D39').Value = '1.294'
$ws.Range('E39').Value = '  -3.81%  '
$ws.Range('D40').Value = '0.6527'
$ws.Range('E40').Value = '  -6.86%  '
$ws.Range('D41').Value = '11.59'
$ws.Range('E41').Value = '  -7.90%  '
$ws.Range('D42').Value = '0.2037'
$ws.Range('E42').Value = '  -8.23%  '
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('E44').Value = '  -7.83%  '
$ws.Range('D45').Value = '13.55'
$ws.Range('E45').Value = '  -5.69%  '
$ws.Range('D46').Value = '2.207'
$ws.Range('E46').Value = '  -5.82%  '
$ws.Range('D47').Value = '1.308'
$ws.Range('E47').Value = '  -5.67%  '
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00000000338'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.08%  '
$ws.Range('D50').Value = '0.06952'
$ws.Range('E50').Value = '  -4.22%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.130'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -8.33%  '
